# LQA_Tester_ProgressTracker.xlsx update
# Refresh the cached QA stats produced by the QAExcelCompiler tool:
#  - DAILY sheet row 8 (01/09): fill in Doni/Eric/Lisa "--" placeholders with
#    real Done/Issues counts now that data exists, and bump John/Mike/Paul
#    totals plus the Pending column.
#  - TOTAL sheet: refresh per-tester rollups (Doni, Eric, Lisa, John, Mike,
#    Paul) and the SUBTOTAL / TOTAL rows.
#  - _DAILY_DATA sheet: refresh underlying per-category rows feeding TOTAL.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# DAILY sheet - row 8 (date 01/09)
# ---------------------------------------------------------------------
$daily = $wb.Worksheets.Item("DAILY")

$daily.Range("F8").Value = 385
$daily.Range("G8").Value = 104
$daily.Range("H8").Value = 3
$daily.Range("I8").Value = 1
$daily.Range("J8").Value = 350
$daily.Range("K8").Value = 5
$daily.Range("N8").Value = 37
$daily.Range("O8").Value = 17
$daily.Range("P8").Value = 253
$daily.Range("Q8").Value = 86
$daily.Range("R8").Value = 254
$daily.Range("X8").Value = 280

# ---------------------------------------------------------------------
# TOTAL sheet - per-tester + subtotal/total rollups
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("TOTAL")

# Row 5 - Doni
$total.Range("B5").Value = 99.7
$total.Range("C5").Value = 100
$total.Range("D5").Value = 385
$total.Range("E5").Value = 104
$total.Range("F5").Value = 276
$total.Range("G5").Value = 5
$total.Range("K5").Value = 104

# Row 6 - Eric
$total.Range("B6").Value = 37.5
$total.Range("C6").Value = 100
$total.Range("D6").Value = 3
$total.Range("E6").Value = 1
$total.Range("F6").Value = 1
$total.Range("G6").Value = 1
$total.Range("K6").Value = 1

# Row 7 - John
$total.Range("D7").Value = 354
$total.Range("E7").Value = 7
$total.Range("F7").Value = 324
$total.Range("G7").Value = 23
$total.Range("K7").Value = 7

# Row 9 - Lisa
$total.Range("B9").Value = 97.40000000000001
$total.Range("C9").Value = 100
$total.Range("D9").Value = 37
$total.Range("E9").Value = 17
$total.Range("F9").Value = 20
$total.Range("K9").Value = 17

# Row 10 - Mike
$total.Range("B10").Value = 96.59999999999999
$total.Range("D10").Value = 253
$total.Range("E10").Value = 86
$total.Range("F10").Value = 167
$total.Range("K10").Value = 86

# Row 11 - Paul
$total.Range("B11").Value = 96.90000000000001
$total.Range("D11").Value = 254
$total.Range("F11").Value = 191

# Row 12 - SUBTOTAL
$total.Range("B12").Value = 95.09999999999999
$total.Range("D12").Value = 1326
$total.Range("E12").Value = 285
$total.Range("F12").Value = 1012
$total.Range("G12").Value = 29
$total.Range("K12").Value = 285

# Row 15 - GRAND TOTAL
$total.Range("B15").Value = 95.09999999999999
$total.Range("D15").Value = 1326
$total.Range("E15").Value = 285
$total.Range("F15").Value = 1012
$total.Range("G15").Value = 29
$total.Range("K15").Value = 285

# ---------------------------------------------------------------------
# _DAILY_DATA sheet - underlying per-category rows
# ---------------------------------------------------------------------
$dailyData = $wb.Worksheets.Item("_DAILY_DATA")

# Row 6 - Doni / Knowledge
$dailyData.Range("D6").Value = 386
$dailyData.Range("E6").Value = 385
$dailyData.Range("F6").Value = 104
$dailyData.Range("G6").Value = 276
$dailyData.Range("H6").Value = 5

# Row 8 - Mike / Region
$dailyData.Range("E8").Value = 253
$dailyData.Range("F8").Value = 86
$dailyData.Range("G8").Value = 167

# Row 9 - Lisa / Region
$dailyData.Range("E9").Value = 37
$dailyData.Range("F9").Value = 17
$dailyData.Range("G9").Value = 20

# Row 10 - Paul / Region
$dailyData.Range("E10").Value = 254
$dailyData.Range("G10").Value = 191

# Row 12 - Eric / Quest
$dailyData.Range("D12").Value = 8
$dailyData.Range("E12").Value = 3
$dailyData.Range("F12").Value = 1
$dailyData.Range("G12").Value = 1
$dailyData.Range("H12").Value = 1

# Row 13 - John / Quest
$dailyData.Range("D13").Value = 8
$dailyData.Range("E13").Value = 7
$dailyData.Range("F13").Value = 2
$dailyData.Range("G13").Value = 1
$dailyData.Range("H13").Value = 4

Write-Host "LQA_Tester_ProgressTracker stats refreshed"
